# Correcao dos casos de uso
# Fixes the UO names used in the MESA_VOTO_ELEICAO insert example statements
# so that each of the 4 example inserts references a different Unidade
# Organica (DEI, DEM, DEEC, FLUC) instead of all of them using 'DEI'.

$d = $word.ActiveDocument

# Avoid Word "smart quotes" autocorrect turning our straight apostrophes
# into curly quotes when we touch text.
try { $word.Options.AutoFormatAsYouTypeReplaceQuotes = $false } catch {}
try { $word.Options.AutoFormatReplaceQuotes = $false } catch {}

function Find-RangeStart($searchText) {
    $r = $d.Content
    $ok = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find text: $searchText"
    }
    return $r.Start
}

# Find the START position of the LAST occurrence of $searchText that begins
# strictly before $beforePos (i.e. the occurrence nearest to, but before,
# beforePos). Search proceeds forward from the top of the document.
function Find-LastRangeStartBefore($searchText, $beforePos) {
    $searchRng = $d.Range(0, 0)
    $lastStart = -1
    $iterations = 0
    while ($true) {
        $iterations = $iterations + 1
        if ($iterations -gt 1000) { break }
        $ok = $searchRng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if (-not $ok) { break }
        if ($searchRng.Start -ge $beforePos) { break }
        $lastStart = $searchRng.Start
        $searchRng.SetRange($searchRng.End, $searchRng.End)
    }
    return $lastStart
}

# Forces a run boundary to exist at the given (collapsed) document position
# without leaving any visible trace, by adding a bookmark there and
# immediately deleting it again.
function Force-RunSplit($pos) {
    $bmName = "zzTempRunSplit"
    if ($d.Bookmarks.Exists($bmName)) {
        $d.Bookmarks($bmName).Delete()
    }
    $d.Bookmarks.Add($bmName, $d.Range($pos, $pos))
    $d.Bookmarks($bmName).Delete()
}

# ---------------------------------------------------------------------
# Block 2 (id = 2): 'DEI' -> 'DEM'
# ---------------------------------------------------------------------
$p = Find-RangeStart "DEI') and id = 2)"
$d.Range($p, $p + 3).Text = "DEM"
Force-RunSplit $p
Force-RunSplit ($p + 3)

# ---------------------------------------------------------------------
# Block 3 (id = 3): 'DEI' -> 'DEEC'
# also splits the preceding "from UNIDADEORGANICA where" run right after
# the "U" (artifact left over from how this text was originally edited).
# ---------------------------------------------------------------------
$p = Find-RangeStart "DEI') and id = 3)"
$uFromPos = Find-LastRangeStartBefore "from UNIDADEORGANICA where" $p
if ($uFromPos -ge 0) {
    Force-RunSplit ($uFromPos + 6)
}
$d.Range($p, $p + 3).Text = "DEEC"

# ---------------------------------------------------------------------
# Block 4 (id = 4): 'DEI' -> 'FLUC'
# the _GoBack bookmark (tracking the last edit location) ends up
# collapsed immediately after the newly typed "FLUC" text.
# ---------------------------------------------------------------------
$p = Find-RangeStart "DEI') and id = 4)"
$d.Range($p, $p + 3).Text = "FLUC"
Force-RunSplit $p
Force-RunSplit ($p + 4)

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $d.Range($p + 4, $p + 4))

Write-Host "Done."
